$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values; dependent formulas (F3, F4, E5, F7) recalc automatically.
$ws.Range("E3").Value = 640
$ws.Range("E4").Value = 640
$ws.Range("E7").Value = 22000

# Update selection to match new active cell in the saved view.
$ws.Range("E8").Select()
